$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write literal text into a cell without Excel's "smart" type
# inference turning date-shaped strings (e.g. "09-08-2018") into serial
# date numbers. We stage the text as a formula result in a scratch cell far
# outside the used range, copy its *value*, and PasteSpecial-Values it into
# the destination - this preserves the destination cell's existing style
# (no quotePrefix / NumberFormat change) while keeping the content as plain
# text (t="s").
function Set-TextValue {
    param($cell, [string]$text)
    $scratch = $ws.Cells.Item(500, 500)
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# --- Fix the incorrect "09-09-2018" dates -> "09-08-2018" ---
Set-TextValue $ws.Cells.Item(12,1) "09-08-2018"
Set-TextValue $ws.Cells.Item(12,2) "09-08-2018"
Set-TextValue $ws.Cells.Item(13,1) "09-08-2018"
Set-TextValue $ws.Cells.Item(13,2) "09-08-2018"
Set-TextValue $ws.Cells.Item(14,1) "09-08-2018"
Set-TextValue $ws.Cells.Item(15,1) "09-08-2018"

# --- Add more prints for the spine box ---
# Row 14 gains a completed date in column B.
Set-TextValue $ws.Cells.Item(14,2) "10-08-2018"

# New row 16: another spine box print.
Set-TextValue $ws.Cells.Item(16,1) "10-08-2018"
Set-TextValue $ws.Cells.Item(16,3) "Spine Box 3lvl Lumbar Spine"
$ws.Cells.Item(16,4).Value = 1
Set-TextValue $ws.Cells.Item(16,5) "HTPLA"
$ws.Cells.Item(16,6).Value = 2
$ws.Cells.Item(16,7).Value = 20
$ws.Cells.Item(16,8).Value = 0.2

# Match the workbook's recorded selection after the edit.
$ws.Range("C17").Select()
